# Applies the 3.9.1.1 worksheet update: adds 2021/2022 data in new
# columns R and S (mirroring the existing D:Q year layout) and moves
# the saved selection to T6, matching the published diff.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Clone the formatting (number format / font / borders / alignment) of the
# existing 2020 column (Q) into the two new year columns before writing values,
# so R/S pick up the same look as the rest of the table.
$ws.Range("Q4:Q34").Copy($ws.Range("R4:R34"))
$ws.Range("Q4:Q34").Copy($ws.Range("S4:S34"))

# Column R = 2021, column S = 2022
$ws.Range("R4").Value = 2021
$ws.Range("S4").Value = 2022
$ws.Range("R5").Value = 0.8
$ws.Range("S5").Value = 0.5
$ws.Range("R6").Value = 0.4
$ws.Range("S6").Value = 0.2
$ws.Range("R7").Value = 1.2
$ws.Range("S7").Value = 0.7
$ws.Range("R8").Value = 0.2
$ws.Range("S8").Value = 0.2
$ws.Range("R9").Value = "-"
$ws.Range("S9").Value = "-"
$ws.Range("R10").Value = 0.4
$ws.Range("S10").Value = 0.4
$ws.Range("R11").Value = 0.6
$ws.Range("S11").Value = 0.5
$ws.Range("R12").Value = 0.8
$ws.Range("S12").Value = 0.3
$ws.Range("R13").Value = 0.5
$ws.Range("S13").Value = 0.6
$ws.Range("R14").Value = 0.4
$ws.Range("S14").Value = 0.7
$ws.Range("R15").Value = "-"
$ws.Range("S15").Value = 0.4
$ws.Range("R16").Value = 0.8
$ws.Range("S16").Value = 1.1000000000000001
$ws.Range("R17").Value = 0.3
$ws.Range("S17").Value = "-"
$ws.Range("R18").Value = 0.7
$ws.Range("S18").Value = "-"
$ws.Range("R19").Value = "-"
$ws.Range("S19").Value = "-"
$ws.Range("R20").Value = 0.5
$ws.Range("S20").Value = 0.4
$ws.Range("R21").Value = 0.1
$ws.Range("S21").Value = 0.4
$ws.Range("R22").Value = 0.8
$ws.Range("S22").Value = 0.4
$ws.Range("R23").Value = 1.1000000000000001
$ws.Range("S23").Value = 0.4
$ws.Range("R24").Value = 1.5
$ws.Range("S24").Value = "-"
$ws.Range("R25").Value = 0.7
$ws.Range("S25").Value = 0.7
$ws.Range("R26").Value = 2.2000000000000002
$ws.Range("S26").Value = 1
$ws.Range("R27").Value = 1
$ws.Range("S27").Value = 0.4
$ws.Range("R28").Value = 3.5
$ws.Range("S28").Value = 1.7
$ws.Range("R29").Value = 0.8
$ws.Range("S29").Value = 0.3
$ws.Range("R30").Value = 0.2
$ws.Range("S30").Value = 0
$ws.Range("R31").Value = 1.6
$ws.Range("S31").Value = 0.6
$ws.Range("R32").Value = 0.3
$ws.Range("S32").Value = "-"
$ws.Range("R33").Value = "-"
$ws.Range("S33").Value = "-"
$ws.Range("R34").Value = 0.6
$ws.Range("S34").Value = "-"

# Restore the saved cursor position recorded in the sheet view
$ws.Range("T6").Select() | Out-Null
